$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.392.08"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.61%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.666.90"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.47%  "
$ws.Range("E4").Value = "  +0.36%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5164"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.88%  "
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.06451"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2575"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.96"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07656"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.341"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.670.01"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.894.04"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5549"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8031"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.66%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.66"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.74%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.403.99"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.56%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.008"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "209.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.421"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.11"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.883"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.47%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.009"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.55"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.54%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.730"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1165"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.004"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.70%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.78"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05232"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.40%  "
$ws.Range("E31").Value = "  -2.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.371"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.220"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.573"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.754"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.381"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.28%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9281"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5714"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.153.20"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +10.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01597"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.97%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.008"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8460"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.646"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.74%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.28"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.72%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.804.43"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.57%  "
$ws.Range("E46").Value = "  -6.39%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4492"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.30%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "55.95"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.50%  "
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.923"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.76%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05113"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.52%  "
